$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 1 (header): shift the old "3/9/2022" header from B1 to C1 and
# add new header cells B1 ("Unnamed: 0.1") and D1:H1 (more dates).
# -----------------------------------------------------------------

# C1 must stay literal text "3/9/2022" (it would otherwise be auto-
# recognised as a date and converted to a serial number), so force a
# text number format before assigning it.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "3/9/2022"

$ws.Range("D1").Value = "14/9/2022"
$ws.Range("E1").Value = "19/9/2022"
$ws.Range("F1").Value = "23/9/2022"
$ws.Range("G1").Value = "27/9/2022"
$ws.Range("H1").Value = "29/9/2022"

# New B1 header cell
$ws.Range("B1").Value = "Unnamed: 0.1"

# Give C1:H1 the same formatting (bold / border / centered) already
# present on B1, reusing the existing style instead of inventing a new one.
$ws.Range("B1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# Column B: item labels (these used to live in column A).
# -----------------------------------------------------------------
$ws.Range("B2").Value = "botas"
$ws.Range("B3").Value = "capacete"
$ws.Range("B4").Value = "colete"
$ws.Range("B5").Value = "luvas"
$ws.Range("B6").Value = "mascara"
$ws.Range("B7").Value = "pa"

# -----------------------------------------------------------------
# Column A: new numeric row-index column (0-5). The old A2:A7 cells
# already carry the bordered/bold style (style "1"), so simply
# overwriting .Value keeps that formatting in place.
# -----------------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# -----------------------------------------------------------------
# Data columns C:H, rows 2-7.
# -----------------------------------------------------------------
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"

$ws.Range("C3").Value = "sim"
$ws.Range("D3").Value = "não"
$ws.Range("E3").Value = "sim"
$ws.Range("F3").Value = "não"
$ws.Range("G3").Value = "não"
$ws.Range("H3").Value = "sim"

$ws.Range("C4").Value = "sim"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("H4").Value = "não"

$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = "sim"

$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "-"
$ws.Range("G6").Value = "-"
$ws.Range("H6").Value = "não"

$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("G7").Value = "-"
$ws.Range("H7").Value = "-"
